# Edit downup 300 names file for clarity
# Visuals for every sweep direction, errors unfixed
#
# Cells A6:A10 used to hold numbers/formulas (a chirp/sweep value ramp) -
# they are relabeled as "bad" (shared string), matching how A11, A57 and
# A62 already flag other rows as bad. The formula chain that used to
# start at A6 (A7=A6+20, A8=A7+20 shared..A31) now starts fresh at A17
# (hard literal 110) with the shared-style ramp continuing from A18
# onward (A18=A17+20, A19=A18+20, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the first five sweep rows as "bad"
$ws.Range("A6:A10").Value = "bad"

# Re-anchor the 110..390 ramp: A17 becomes a literal 110, and the
# remaining rows continue the +20 pattern from there.
$ws.Range("A17").Formula = "=110"
$ws.Range("A18").Formula = "=A17+20"

# Move the active selection to A10 (this also clears any scrolled
# viewport, so the sheet view opens at the top-left).
[void]$ws.Range("A10").Select()
